# Auto update Excel log
# Appends newly-logged sensor events to the relevant sheets of the
# SeniorConnect master log workbook.

$wb = $excel.ActiveWorkbook

function Add-LogRows($SheetName, $StartRow, $Rows) {
    $ws = $wb.Worksheets.Item($SheetName)
    $r = $StartRow
    foreach ($row in $Rows) {
        # Column A holds a literal "YYYY-MM-DD" string; force text so it
        # is not reinterpreted as a date serial number.
        $dateCell = $ws.Cells.Item($r, 1)
        $dateCell.NumberFormat = "@"
        $dateCell.Value = $row[0]

        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $ws.Cells.Item($r, 4).Value = $row[3]
        $ws.Cells.Item($r, 5).Value = $row[4]
        $ws.Cells.Item($r, 6).Value = $row[5]

        $r = $r + 1
    }
}

# --- PIR sheet: two new "Out of Bed" events in the Bedroom ---
$pirRows = @(
    @("2026-02-01", "14:46:17", "14:00", "Bedroom", "Out of Bed", "Empty"),
    @("2026-02-01", "14:46:21", "14:00", "Bedroom", "Out of Bed", "Empty")
)
Add-LogRows "PIR" 2 $pirRows

# --- Proximity sheet: four new door ENTER/EXIT events ---
$proximityRows = @(
    @("2026-02-01", "14:46:01", "14:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-02-01", "14:46:13", "14:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-02-01", "14:46:17", "14:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door"),
    @("2026-02-01", "14:46:47", "14:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door")
)
Add-LogRows "Proximity" 48 $proximityRows

# --- Camera sheet: two new image events ---
$cameraRows = @(
    @("2026-02-01", "14:46:03", "14:00", "Living Room Main Door", "Image Received", "Active"),
    @("2026-02-01", "14:46:17", "14:00", "Living Room Main Door", "Image Captured", "Active")
)
Add-LogRows "Camera" 33 $cameraRows

# --- mmWave(BR) / mmWave(HR) sheets: same two new Bedroom readings
#     (numeric Value column = 0, not text) ---
$mmWaveRows = @(
    @("2026-02-01", "14:46:18", "14:00", "Bedroom", 0, "Empty"),
    @("2026-02-01", "14:46:21", "14:00", "Bedroom", 0, "Empty")
)
Add-LogRows "mmWave(BR)" 2 $mmWaveRows
Add-LogRows "mmWave(HR)" 2 $mmWaveRows
